$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New tag value in A2 (typed first, per shared-string ordering).
$ws.Range("A2").Value = "PT_62"

# Insert "NPW_array_PT_2084" above the old B3 (NPW_array_PT_61). Push column-B
# values down one row at a time, bottom-up so nothing is clobbered before it's
# read; this leaves the other columns (C3/F3/I3/L3 marker cells) untouched,
# matching Excel's single-column Range.Insert behaviour.
for ($r = 23; $r -ge 3; $r--) {
    $v = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 2).Value = $v
}
$ws.Cells.Item(3, 2).Value = "NPW_array_PT_2084"

# NPW_array_PT_4025 is now at row 19; insert "NPW_array_PT_4013" above it,
# shifting B19:B24 down to B20:B25.
for ($r = 24; $r -ge 19; $r--) {
    $v = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 2).Value = $v
}
$ws.Cells.Item(19, 2).Value = "NPW_array_PT_4013"

# NPW_array_PT_5000A is now at row 25; insert "NPW_array_PT_5013" above it,
# shifting B25 down to B26.
$ws.Cells.Item(26, 2).Value = $ws.Cells.Item(25, 2).Value2
$ws.Cells.Item(25, 2).Value = "NPW_array_PT_5013"

# Update the sheet's selection to match the final state.
$ws.Range("C25").Select()
